$wb = $excel.ActiveWorkbook

# --- About sheet: bump the "last updated" date ---
$about = $wb.Worksheets.Item("About")
$about.Range("C1").Value = "3/28/2024"

# --- RAF-capacity sheet: update hydrogen RAF values to full credit ---
$capacity = $wb.Worksheets.Item("RAF-capacity")
$capacity.Range("B24").Value = 1
$capacity.Range("B25").Value = 1

# --- Resize column A on RAF-capacity to fit the shorter labels ---
$capacity.Columns.Item(1).ColumnWidth = 28.166666666666668

# --- Switch the active/visible tab from RAF-generation to RAF-capacity ---
$capacity.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 80
$capacity.Range("B25").Select() | Out-Null
